$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column M ("Gas Choice ID"), shifting the
# existing Gas Choice ID / Gas Rate Code / Gas Usage columns one to the
# right (M->N, N->O, O->P) and growing the sheet's used range to A1:P2.
$ws.Range("M1").EntireColumn.Insert()

# Give the freshly inserted header cell the same formatting as its
# neighboring header cells (bold, bordered, centered) by copying the
# format from the adjacent header cell, then set its text.
$ws.Range("N1").Copy()
$ws.Range("M1").PasteSpecial(-4122)
$ws.Range("M1").Value = "Gas Supplier"

# New data cell for row 2 under the inserted "Gas Supplier" column.
$ws.Range("M2").Value = "N/A"

$excel.CutCopyMode = 0
